$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9911111111111112
$ws.Range("C2").Value = 0.9908032995990874
$ws.Range("D2").Value = 0.9911162109012718
$ws.Range("E2").Value = 0.9901065303769977
$ws.Range("F2").Value = 28.6
$ws.Range("G2").Value = 2469.582492828369

$ws.Range("B3").Value = 0.9901851851851852
$ws.Range("C3").Value = 0.9898819612599488
$ws.Range("D3").Value = 0.9901915809835945
$ws.Range("E3").Value = 0.9890765443407569
$ws.Range("F3").Value = 28.3
$ws.Range("G3").Value = 1474.029290914536

$ws.Range("C4").Value = 0.9895197400624811
$ws.Range("D4").Value = 0.9898137747893395
$ws.Range("E4").Value = 0.9886636029131105
$ws.Range("F4").Value = 12.2
$ws.Range("G4").Value = 1558.9833984375

$ws.Range("A5").Value = "densenet121"
$ws.Range("B5").Value = 0.9896296296296296
$ws.Range("C5").Value = 0.9893331952174449
$ws.Range("D5").Value = 0.9896173866341067
$ws.Range("E5").Value = 0.9884574425448591
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 1685.513358354568

$ws.Range("A6").Value = "resnet101"
$ws.Range("B6").Value = 0.9890740740740741
$ws.Range("C6").Value = 0.9887031337958833
$ws.Range("D6").Value = 0.9890730406388604
$ws.Range("E6").Value = 0.9878390912526194
$ws.Range("F6").Value = 44.5
$ws.Range("G6").Value = 1928.099274635315

$ws.Range("B7").Value = 0.9875925925925926
$ws.Range("C7").Value = 0.9870486448163949
$ws.Range("D7").Value = 0.9875855830309455
$ws.Range("E7").Value = 0.9861907342804144
$ws.Range("F7").Value = 5.3
$ws.Range("G7").Value = 1087.654934644699

$ws.Range("A8").Value = "vit_b_16"
$ws.Range("B8").Value = 0.9872222222222222
$ws.Range("C8").Value = 0.9869984774436796
$ws.Range("D8").Value = 0.9872227458213874
$ws.Range("E8").Value = 0.9857794943820225
$ws.Range("F8").Value = 86.59999999999999
$ws.Range("G8").Value = 3402.739871740341

$ws.Range("A9").Value = "resnet50"
$ws.Range("B9").Value = 0.987037037037037
$ws.Range("C9").Value = 0.9865022853476677
$ws.Range("D9").Value = 0.987038332097363
$ws.Range("E9").Value = 0.9855716379626237
$ws.Range("F9").Value = 25.6
$ws.Range("G9").Value = 1629.672913789749
